$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.349.32"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.643.44"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.14"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.58"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "3.110.11"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "59.349.52"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.87"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000137"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.670.08"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "342.72"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.45"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.63"
$ws.Range("E20").Value = "  +2.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.40"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.67"
$ws.Range("E23").Value = "  +3.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.417"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.21"
$ws.Range("E27").Value = "  +1.60%  "
$ws.Range("D28").Value = "0.0₃0803"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.40"
$ws.Range("E30").Value = "  -4.13%  "
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.04"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.77"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.16"
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.872"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.864"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.49"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.54"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.65"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0977"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.601"
$ws.Range("E43").Value = "  -3.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "270.25"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.47"
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0536"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").Value = "2.037.80"
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.77"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0229"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.78"
$ws.Range("E51").Value = "  -1.40%  "
